$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("StageCoordinate")

# --- Row 2 ---
$ws.Range("C2").Value = 2000
$ws.Range("D2").Value = 2000

# --- Row 3 ---
$ws.Range("C3").Value = 3000
$ws.Range("D3").Value = 3000
$ws.Range("C5").Copy()
$ws.Range("D3").PasteSpecial(-4122)  # xlPasteFormats

# --- Row 4 (stage text swaps with row 6's old text; values become 3000) ---
$ws.Range("B4").Value = "-1,-1,-1,-1,4,-1,-1,-1,-1,-1,-1,-1,0,-1,-1,1,-1,-1,-1,3,-1,0,-1,-1,-1"
$ws.Range("C4").Value = 3000
$ws.Range("D4").Value = 3000
$ws.Range("C5").Copy()
$ws.Range("C4").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("C5").Copy()
$ws.Range("D4").PasteSpecial(-4122)  # xlPasteFormats

# --- Row 5 (new, previously-unused stage text; values become 4000) ---
$ws.Range("B5").Value = "-1,-1,-1,-1,4,-1,-1,0,-1,3,1,-1,-1,2,-1,-1,-1,-1,-1,-1,-1,-1,-1,-1,-1"
$ws.Range("C5").Value = 4000
$ws.Range("D5").Value = 4000

# --- Row 6 (stage text swaps with row 5's old text; values become 4000) ---
$ws.Range("B6").Value = "-1,-1,-1,-1,3,-1,-1,0,-1,-1,4,-1,-1,2,-1,-1,1,-1,-1,-1,-1,-1,-1,0,-1"
$ws.Range("C6").Value = 4000
$ws.Range("D6").Value = 4000

# --- Selection ---
$excel.CutCopyMode = $false
$ws.Activate()
$ws.Range("D6").Select()
